# Add new columns I (I0) and J (IF) to the worksheet, mirroring column H's
# header style, and populate the data rows 2-60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered alignment) from
# the existing H1 header cell so the new headers look consistent.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data values --------------------------------------------------------
$iValues = @(8,7,10,8,8,6,7,8,8,8,9,7,7,5,8,7,7,8,9,7,7,8,8,7,8,6,8,8,7,8,7,8,7,7,5,8,8,8,8,7,8,8,8,8,8,8,8,8,8,8,8,8,9,4,3,3,6,4,3)
$jValues = @(8,7,10,8,8,7,7,8,8,8,9,7,7,5,8,7,8,8,9,7,7,8,8,8,8,8,8,8,7,8,8,9,8,7,6,8,8,8,8,8,8,8,8,8,8,8,8,8,8,8,8,8,9,4,3,3,6,4,3)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]   # column I
    $ws.Cells.Item($row, 10).Value = $jValues[$i]  # column J
}
